$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one row to accommodate the new Adelson/Africana (Clarke) rows (122 -> 123 data rows)
$ws.Rows.Item(14).Insert()

# Rewrite the full data block (rows 3-123, columns A-E) to reflect the updated
# Adelson/Africana Library naming plus the resulting row reflow.
$data = New-Object 'object[,]' 121,5
$data[0,0] = 'Adelson Library'
$data[0,4] = 'Adelson Library > Main Collection'
$data[1,0] = 'Adelson Library Reference (Non-Circulating)'
$data[1,4] = 'Adelson Library > Reference'
$data[2,0] = 'Museum of Vertebrates (Non-Circulating)'
$data[2,4] = 'Adelson Library > Museum of Vertebrates'
$data[3,0] = 'Macaulay Library (Non-Circulating)'
$data[3,4] = 'Adelson Library > Macaulay Library'
$data[4,0] = 'Clarke Africana Library A/V'
$data[4,4] = 'Africana Library > A/V'
$data[5,0] = 'Clarke Africana Library New and Noteworthy'
$data[5,4] = 'Africana Library > New Books Shelf'
$data[6,0] = 'Clarke Africana Library'
$data[6,3] = 'new book'
$data[6,4] = 'Africana Library > New Books Shelf'
$data[7,0] = 'Clarke Africana Library'
$data[7,4] = 'Africana Library > Main Collection'
$data[8,0] = 'Clarke Africana Library Reference ( Non-Circulating)'
$data[8,4] = 'Africana Library > Reference'
$data[9,0] = 'Clarke Africana Library Reserve'
$data[9,4] = 'Africana Library > Reserve'
$data[10,0] = 'Clarke Africana Library Permanent Reserve'
$data[10,4] = 'Africana Library > Reserve'
$data[11,0] = 'Bailey Hortorium (ask at Mann Library Circulation)'
$data[11,4] = 'Bailey Hortorium > Main Collection'
$data[12,0] = 'Bailey Hortorium Reference (Non-Circulating)'
$data[12,4] = 'Bailey Hortorium > Reference'
$data[13,0] = 'CISER Data Archive'
$data[13,4] = 'CISER Data Archive'
$data[14,0] = 'Fine Arts Circulation '
$data[14,4] = 'Fine Arts Library > Circulation'
$data[15,0] = 'Fine Arts Library (Rand Hall)'
$data[15,4] = 'Fine Arts Library > Main Collection'
$data[16,0] = 'Fine Arts Library Reference (Non-Circulating)'
$data[16,2] = 'Artist'
$data[16,4] = 'Fine Arts Library > Artists'' Books'
$data[17,0] = 'Fine Arts Library Reference (Non-Circulating)'
$data[17,4] = 'Fine Arts Library > Reference'
$data[18,0] = 'Fine Arts Library Permanent Reserve'
$data[18,4] = 'Fine Arts Library > Reserve'
$data[19,0] = 'Fine Arts Course Reserve (Ask at Circulation)'
$data[19,4] = 'Fine Arts Library > Reserve'
$data[20,0] = 'Fine Arts Library Reserve'
$data[20,4] = 'Fine Arts Library > Reserve'
$data[21,0] = 'ILR Library (Ives Hall)'
$data[21,4] = 'ILR Library > Main Collection'
$data[22,0] = 'ILR Library Reference (Non-Circulating)'
$data[22,2] = 'Labor Law'
$data[22,4] = 'ILR Library > Labor Law Reference'
$data[23,0] = 'ILR Library Reference (Non-Circulating)'
$data[23,4] = 'ILR Library > Reference'
$data[24,0] = 'ILR Permanent Reserve '
$data[24,4] = 'ILR Library > Reserve'
$data[25,0] = 'ILR Library Reserve '
$data[25,4] = 'ILR Library > Reserve'
$data[26,0] = 'ILR Multi-Copy Storage'
$data[26,4] = 'ILR Library > Main Collection'
$data[27,0] = 'ILR Library Kheel Center'
$data[27,4] = 'ILR Library Kheel Center > Main Collection'
$data[28,0] = 'ILR Library Kheel Center (Non-Circulating)'
$data[28,4] = 'ILR Library Kheel Center > Main Collection'
$data[29,0] = 'ILR Library Kheel Center (Request in advance)'
$data[29,4] = 'ILR Library Kheel Center > Request in advance'
$data[30,0] = 'ILR Library Kheel Center Reference '
$data[30,4] = 'ILR Library Kheel Center > Reference'
$data[31,0] = 'ILR Library Kheel Center Reference (Non-Circulating)'
$data[31,4] = 'ILR Library Kheel Center > Reference'
$data[32,0] = 'Asia Reserve, Severinghouse Reading Rm., Kroch Library'
$data[32,4] = 'Kroch Library Asia > Reserve'
$data[33,0] = 'Kroch Library Asia'
$data[33,4] = 'Kroch Library Asia > Main Collection'
$data[34,0] = 'Kroch Library Asia Reference (Non-Circulating)'
$data[34,4] = 'Kroch Library Asia > Reference'
$data[35,0] = 'RMC Technical Services'
$data[35,1] = 'X'
$data[35,4] = 'Kroch Library Rare & Manuscripts > Technical Services'
$data[36,0] = 'Kroch Library Rare & Manuscripts (Non-Circulating)'
$data[36,4] = 'Kroch Library Rare & Manuscripts > Main Collection'
$data[37,0] = 'Kroch Library Rare & Manuscripts (Request in advance)'
$data[37,4] = 'Kroch Library Rare & Manuscripts > Request in Advance'
$data[38,0] = 'Kroch Library Rare & Manuscripts Reference (Non-Circulating)'
$data[38,4] = 'Kroch Library Rare & Manuscripts > Reference'
$data[39,0] = 'Law Library (Myron Taylor Hall)'
$data[39,4] = 'Law Library > Main Collection'
$data[40,0] = 'Law Library (Myron Taylor Hall) Rare Books'
$data[40,4] = 'Law Library > Rare Books'
$data[41,0] = 'Law Library Rare--Request in advance at Law Circulation Desk'
$data[41,4] = 'Law Library > Request in Advance'
$data[42,0] = 'Law Library Reference (Non-Circulating)'
$data[42,4] = 'Law Library > Reference'
$data[43,0] = 'Law Library Reserve'
$data[43,4] = 'Law Library > Reserve'
$data[44,0] = 'Law Library Technical Services '
$data[44,4] = 'Law Library > Technical Services'
$data[45,0] = 'Legal Aid Clinic'
$data[45,4] = 'Law Library > Legal Aid Clinic'
$data[46,0] = 'Library Annex'
$data[46,4] = 'Library Annex'
$data[47,0] = 'Mann Circulation'
$data[47,4] = 'Mann Library > Circulation'
$data[48,0] = 'Mann Library'
$data[48,2] = 'Ellis'
$data[48,4] = 'Mann Library > Ellis Collection'
$data[49,0] = 'Mann Library'
$data[49,2] = 'Curriculum Material'
$data[49,4] = 'Mann Library > Curriculum Materials Collection'
$data[50,0] = 'Mann Library'
$data[50,4] = 'Mann Library > Main Collection'
$data[51,0] = 'Mann Library Collection Development (Non-Circulating)'
$data[51,4] = 'Mann Library > Collection Development'
$data[52,0] = 'Mann Library New Book Shelf'
$data[52,4] = 'Mann Library > New Book Shelf'
$data[53,0] = 'Mann Library Reference (Non-Circulating)'
$data[53,4] = 'Mann Library > Reference'
$data[54,0] = 'Mann Library Reserve'
$data[54,4] = 'Mann Library > Reserve'
$data[55,0] = 'Mann Library Special Collections (Non-Circulating)'
$data[55,4] = 'Mann Library > Special Collections'
$data[56,0] = 'Mann Special Collections (Request in advance)'
$data[56,4] = 'Mann Library > Special Collections'
$data[57,0] = 'Mann Serials'
$data[57,4] = 'Mann Library > Serials'
$data[58,0] = 'Mathematics Library (Circulation Desk)'
$data[58,4] = 'Mathematics Library > Circulation Desk'
$data[59,0] = 'Mathematics Library (Malott Hall)'
$data[59,4] = 'Mathematics Library > Main Collection'
$data[60,0] = 'Mathematics Library Locked Press'
$data[60,4] = 'Mathematics Library > Locked Press'
$data[61,0] = 'Mathematics Library Reference (Non-Circulating)'
$data[61,4] = 'Mathematics Library > Reference'
$data[62,0] = 'Mathematics Library Reserve'
$data[62,4] = 'Mathematics Library > Reserve'
$data[63,0] = 'Music Library (Lincoln Hall)'
$data[63,4] = 'Music Library > Main Collection'
$data[64,0] = 'Cox Library of Music (Lincoln Hall)'
$data[64,4] = 'Music Library > Main Collection'
$data[65,0] = 'Music Circulation'
$data[65,4] = 'Music Library > Circulation'
$data[66,0] = 'Music Library A/V (Non-Circulating)'
$data[66,4] = 'Music Library > A/V'
$data[67,0] = 'Music Library Locked Press (Reference Desk)'
$data[67,4] = 'Music Library > Reference'
$data[68,0] = 'Music Library Reference (Non-Circulating) '
$data[68,4] = 'Music Library > Reference'
$data[69,0] = 'Music Library Reserve '
$data[69,4] = 'Music Library > Reserve'
$data[70,0] = 'Nestle Library Permanent Reserve'
$data[70,4] = 'Nestle Library > Permanent Reserve'
$data[71,0] = 'Nestle Library Reserve '
$data[71,4] = 'Nestle Library > Reserve'
$data[72,0] = 'Olin Library'
$data[72,2] = 'New & Noteworthy'
$data[72,4] = 'Olin Library > New & Noteworthy Books Shelf'
$data[73,0] = 'Olin Library'
$data[73,4] = 'Olin Library > Main Collection'
$data[74,0] = 'DCAPS (106G Olin)'
$data[74,4] = 'Olin Library > DCAPS (106G Olin)'
$data[75,0] = 'Map Storage (Request in Advance at Map Room, Olin Library)'
$data[75,4] = 'Olin Library > Map Storage'
$data[76,0] = 'Olin Library Maps (Non-Circulating)'
$data[76,4] = 'Olin Library > Maps'
$data[77,0] = 'Olin Library Graduate Study Room 501, Request at Circulation'
$data[77,4] = 'Olin Library > Graduate Study Room 501'
$data[78,0] = 'Olin Library Reference (Non-Circulating)'
$data[78,4] = 'Olin Library > Reference'
$data[79,0] = 'Olin Library Reserve'
$data[79,4] = 'Olin Library > Reserve'
$data[80,0] = 'Olin Library Room 301 (Non-Circulating)'
$data[80,4] = 'Olin Library > Room 301'
$data[81,0] = 'Olin Library Room 305 (Non-Circulating)'
$data[81,4] = 'Olin Library > Room 305'
$data[82,0] = 'Olin Library Room 401 (Non-Circulating)'
$data[82,4] = 'Olin Library > Room 401'
$data[83,0] = 'Olin Library Room 404 (Non-Circulating)'
$data[83,4] = 'Olin Library > Room 404'
$data[84,0] = 'Olin Library Room 405 (Non-Circulating)'
$data[84,4] = 'Olin Library > Room 405'
$data[85,0] = 'Olin Library Room 602 (Non-Circulating)'
$data[85,4] = 'Olin Library > Room 602'
$data[86,0] = 'Olin Library Room 604-605 (Non-Circulating)'
$data[86,4] = 'Olin Library > Room 604-605'
$data[87,0] = 'Olin Library Room 303'
$data[87,4] = 'Olin Library > Room 303'
$data[88,0] = 'Olin Library Room 403'
$data[88,4] = 'Olin Library > Room 403'
$data[89,0] = 'Olin Library Room 603'
$data[89,4] = 'Olin Library > Room 603'
$data[90,0] = 'OKU Processing'
$data[90,1] = 'X'
$data[90,4] = 'Olin Library > OKU Processing'
$data[91,0] = 'Library Tech Services'
$data[91,1] = 'X'
$data[91,4] = 'Olin Library > Technical Services'
$data[92,0] = 'Library Technical Services Review Shelves'
$data[92,1] = 'X'
$data[92,4] = 'Olin Library > Technical Services Review Shelves'
$data[93,0] = 'LTS E-Resources and Serials Management'
$data[93,1] = 'X'
$data[93,4] = 'Olin Library > Technical Services'
$data[94,0] = 'Preservation Department (B32 Olin)'
$data[94,4] = 'Olin Library > Preservation Department'
$data[95,0] = 'Request at Olin Circulation Desk'
$data[95,4] = 'Olin Library > Circulation Desk'
$data[96,0] = 'Sage Hall Management Library Reference (Non-Circulating)'
$data[96,4] = 'Sage Hall Management Library > Reference'
$data[97,0] = 'Sage Hall Management Library Reserve'
$data[97,4] = 'Sage Hall Management Library > Reserve'
$data[98,0] = 'Spacecraft Planetary Imaging Facility (Non-Circulating)'
$data[98,4] = 'Space Sciences Building'
$data[99,0] = 'Uris Library'
$data[99,4] = 'Uris Library > Main Collection'
$data[100,0] = 'Uris Circulation'
$data[100,4] = 'Uris Library > Circulation'
$data[101,0] = 'Uris Library Asia A/V'
$data[101,2] = '(SEA)'
$data[101,4] = 'Uris Library > Asia A/V > Southeast Asia'
$data[102,0] = 'Uris Library Asia A/V'
$data[102,2] = '(SA)'
$data[102,4] = 'Uris Library > Asia A/V > South Asia'
$data[103,0] = 'Uris Library Asia A/V'
$data[103,2] = '(KOR)'
$data[103,4] = 'Uris Library > Asia A/V > Korea'
$data[104,0] = 'Uris Library Asia A/V'
$data[104,2] = '(Jpn)'
$data[104,4] = 'Uris Library > Asia A/V > Japan'
$data[105,0] = 'Uris Library Asia A/V'
$data[105,2] = '(Chi)'
$data[105,4] = 'Uris Library > Asia A/V > China'
$data[106,0] = 'Uris Library Asia A/V'
$data[106,4] = 'Uris Library > Asia A/V'
$data[107,0] = 'Uris Library Reference (Non-Circulating)'
$data[107,4] = 'Uris Library > Reference'
$data[108,0] = 'Uris Library Reserve'
$data[108,4] = 'Uris Library > Reserve'
$data[109,0] = 'Uris Library Reserve Willis Room '
$data[109,4] = 'Uris Library > Reserve Willis Room '
$data[110,0] = 'Veterinary Library (Schurman Hall)'
$data[110,4] = 'Veterinary Library > Main Collection'
$data[111,0] = 'Veterinary Library Core Resource (5 hour loan)'
$data[111,4] = 'Veterinary Library > Core Textbooks'
$data[112,0] = 'Veterinary Library Rare Books (Non-Circulating)'
$data[112,4] = 'Veterinary Library > Rare Books'
$data[113,0] = 'Veterinary Library Reference (Non-Circulating)'
$data[113,4] = 'Veterinary Library > Reference'
$data[114,0] = 'Veterinary Library Reserve '
$data[114,4] = 'Veterinary Library > Reserve '
$data[115,0] = 'Anatomic Pathology Collection (Departmental use only)'
$data[115,1] = 'X'
$data[115,4] = 'Veterinary Library > Anatomic Pathology Collection'
$data[116,0] = 'Equine Farm Animal Collection (Departmental use only)'
$data[116,1] = 'X'
$data[116,4] = 'Veterinary Library > Equine Farm Animal Collection'
$data[117,0] = 'Center for Animal Resources and Education (Dept. use only)'
$data[117,1] = 'X'
$data[117,4] = 'XXX'
$data[118,0] = 'Companion Animal Hospital Collection (Departmental use only)'
$data[118,1] = 'X'
$data[118,4] = 'XXX'
$data[119,0] = 'Clinical Ophthalmology Collection (Departmental use only)'
$data[119,1] = 'X'
$data[119,4] = 'XXX'
$data[120,0] = 'Networked Resource'
$data[120,1] = 'X'
$data[120,4] = 'XXX'

$ws.Range("A3:E123").Value2 = $data

$ws.Range("A12").Select()
